$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update marketCapUsd (column F) values for rows 2-101 (excluding the swapped rows 88-89, handled below)
$marketCapUpdates = @{
    2 = "1645600371530.6442755914265066"
    3 = "221238742768.5913530967335473"
    4 = "144656039999.9366479492830614"
    5 = "124997297363.8779204833653026"
    6 = "87163354463.4157238705464440"
    7 = "64106576736.5184495031453656"
    8 = "60326188982.8522269574726117"
    9 = "25323753239.1835317661563318"
    10 = "23771977841.4471850239419746"
    11 = "21875136033.2702502834964841"
    12 = "17108030022.0500621826818085"
    13 = "10622166943.4478589432961556"
    14 = "9118365895.6041472765326977"
    15 = "8908742037.0967770360407129"
    16 = "8861972999.7187600514796678"
    17 = "8236068214.6058648044301618"
    18 = "8167767101.3844922988942224"
    19 = "7487771074.9913127067659214"
    20 = "7319371445.2538029592242559"
    21 = "7291815763.0727263280639646"
    22 = "6469254257.0149880701896458"
    23 = "6211954278.5391776717943300"
    24 = "6144102402.9310631778808316"
    25 = "6104196066.7293727320696468"
    26 = "6023188746.3479994497551575"
    27 = "5577893743.8175652400000000"
    28 = "5406233793.7494721118674033"
    29 = "5364422326.3378411111389745"
    30 = "5245728803.7280329393707867"
    31 = "4298824719.0012374183297740"
    32 = "3993800282.9340974800817414"
    33 = "3769326792.3955632119224379"
    34 = "3154710197.9810056551877096"
    35 = "3097008436.6592812081353600"
    36 = "3029773532.5168219759114819"
    37 = "2891620175.4399099420000000"
    38 = "2790157658.1501046850958512"
    39 = "2699961505.2057555335563040"
    40 = "2584061240.9943523653083949"
    41 = "2581685391.1439520595162100"
    42 = "2549908349.0321047639046086"
    43 = "2531820394.5145716138260667"
    44 = "2516755433.8775290546824364"
    45 = "2109448878.9105546162727262"
    46 = "2017932366.9006062675018450"
    47 = "1972372076.7056276660095994"
    48 = "1968864590.6431538780648340"
    49 = "1931516772.0549792309375000"
    50 = "1914050754.4524487211698307"
    51 = "1861481627.6244087027573087"
    52 = "1801854633.9552019550775493"
    53 = "1767555384.7792514227107219"
    54 = "1667781804.2086878000182569"
    55 = "1654706339.0678902770221032"
    56 = "1572698970.0981830710708801"
    57 = "1541312353.2150975275048328"
    58 = "1467334025.4275592598641476"
    59 = "1410827960.6605581120000000"
    60 = "1369938502.8672657170340568"
    61 = "1356745023.4707309026682201"
    62 = "1283915640.4706094607607447"
    63 = "1246676762.6033562908945136"
    64 = "1191234770.1155109600000000"
    65 = "1099987973.7117743500000000"
    66 = "1097148630.1213933849797655"
    67 = "1070156987.8680359581709838"
    68 = "976331197.3122773245312771"
    69 = "961915715.5932057785932830"
    70 = "951874736.5640840923995073"
    71 = "944028010.9010396970519914"
    72 = "898582459.0946837021955001"
    73 = "885191982.2847630618795874"
    74 = "882390168.3221567905715871"
    75 = "873887205.0400928502751168"
    76 = "868354716.6944952243536303"
    77 = "851562953.7912270973863332"
    78 = "840759571.2912169595663004"
    79 = "836039972.3196338000000000"
    80 = "830198277.1416908051489310"
    81 = "806973804.9338775850067293"
    82 = "801516669.7799929918057961"
    83 = "795127890.3225529802920000"
    84 = "788485386.8367231244493098"
    85 = "767629346.5031013806052952"
    86 = "696775670.0897274689012979"
    87 = "684048671.4149622103839897"
    90 = "670010329.3811026259130289"
    91 = "663315970.4451091582092000"
    92 = "656440612.0699848713241346"
    93 = "643333898.4821957136702312"
    94 = "642893675.9930117469782132"
    95 = "618550480.8274129470079088"
    96 = "612140479.4307130801030151"
    97 = "603786260.1301708746558791"
    98 = "594222462.9006889619505658"
    99 = "579262544.1824133479025006"
    100 = "572038782.6309731994196095"
    101 = "565351361.5889217448536238"
}
foreach ($row in $marketCapUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 6)
    $cell.NumberFormat = "@"
    $cell.Value = $marketCapUpdates[$row]
}

# Row 88 and 89 swap their id/name/supply/explorer content (rank in column C is unchanged),
# and both rows also get updated marketCapUsd (column F) values.
$ws.Cells.Item(88, 1).Value = "gala"
$ws.Cells.Item(88, 2).Value = "Gala"
$ws.Cells.Item(88, 4).NumberFormat = "@"
$ws.Cells.Item(88, 4).Value = "43744189534.5928200000000000"
$ws.Cells.Item(88, 5).Value = "https://ethplorer.io/es/address/0x15d4c048f83bd7e37d49ea4c83a07267ec4203da#chart=candlestick"
$ws.Cells.Item(88, 6).NumberFormat = "@"
$ws.Cells.Item(88, 6).Value = "675934600.2598317924475358"

$ws.Cells.Item(89, 1).Value = "tezos"
$ws.Cells.Item(89, 2).Value = "Tezos"
$ws.Cells.Item(89, 4).NumberFormat = "@"
$ws.Cells.Item(89, 4).Value = "1036827531.6780720000000000"
$ws.Cells.Item(89, 5).Value = "https://tzkt.io/"
$ws.Cells.Item(89, 6).NumberFormat = "@"
$ws.Cells.Item(89, 6).Value = "673831428.5401824338962507"
